$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change test_number (column A) to 4 for sample rows 2-11
$ws.Range("A2:A11").Value = 4

# Update answers_4 (column E) for row 3 from 5 to 2
$ws.Range("E3").Value = 2

# Update the active selection to I10 (matches recorded cursor position)
$ws.Range("I10").Select()
